$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games), columns B:G for rows 2-10.
# F (Win) is unchanged; G = B + C + D + E (sum column).

$data = @{
    2  = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    3  = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 3.537761648806719;  E = 0.4942365360607697; G = 8.974608811992548 }
    4  = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    5  = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    6  = @{ B = 0.6606524410359556;   C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 2.960089034096801 }
    7  = @{ B = 0.00001292064567892659; C = 10.34677158129881; D = 0.7527432677738641; E = 10.19245300693656;  G = 21.29198077665491 }
    8  = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    9  = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    10 = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 3.537761648806719;  E = 0.4942365360607697; G = 8.974608811992548 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
